$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Hortaliza, Zapallo - variedad "Paine") is inserted
# as row 84; every existing record from row 84 down shifts one row lower
# (old 84 -> 85, ..., old 98 -> 99). Insert() on the row shifts cells down
# and carries the formatting of the row above, matching Excel's native
# "Insert" behaviour (this also carries the date style onto the new D84).
$ws.Rows("84").Insert()

$ws.Range("A84").Value = 7
$ws.Range("B84").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C84").Value = "Ñuble"
$ws.Range("D84").Value = 44505
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 100112045
$ws.Range("G84").Value = "Zapallo"
$ws.Range("H84").Value = "Paine"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 600
$ws.Range("K84").Value = 180
$ws.Range("L84").Value = 200
$ws.Range("M84").Value = 190
$ws.Range("N84").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O84").Value = "Región de O'Higgins"
$ws.Range("P84").Value = 190
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"
